# Updated Input files and Indexes
#
# The "StudyFilesTab" query stored in cell B5 of Sheet1 is rewritten so
# that every reference to `sf.original_file_size` becomes `sf.file_size`
# (the underlying data model renamed that column). Re-saving this edited
# string naturally moves it to the end of the shared-string table, which
# is why B2/B3/B4's shared-string indexes also shift down - that is an
# automatic side effect of the text edit, not a separate change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B5")
$oldText = $cell.Value()
$newText = $oldText.Replace("original_file_size", "file_size")
$cell.Value = $newText

# Re-assert the row's height: editing the wrapped text re-triggers
# autosizing, so pin it back to its original value.
$ws.Rows.Item(5).RowHeight = 409.6

# Reflect the window/selection state at the point the edit was made:
# scrolled down one row (top-left cell A4) with C5 selected.
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C5").Select()
